# Refresh the live crypto price/volume snapshot (GitHub Actions cron update).
# Values are written as plain text, matching the source sheet's inline-string cells
# (price/volume columns are formatted text, not numeric, in this workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.512.47'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '3.778.26'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.38'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000248'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('D14').Value = '4.412.94'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '3.790.60'
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').Value = '67.568.97'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.32'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.02'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '459.23'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.69'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.693'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('E23').Value = '  -4.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.45'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.92'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').Value = '3.926.46'
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +2.90%  '
$ws.Range('E31').Value = '  -6.15%  '
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.03'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.95'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0988'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.24'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.987'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '47.36'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('B44').Value = 'Arweave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.35'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.296'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.44'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.30'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.35'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '26.90'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '389.51'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.31%  '
